# BALANCING.xlsx — extend the compounding-growth table from 50 rows (2-51)
# to 100 rows (2-101) and update the view/selection to match.
#
# Column A holds the running index (A(r) = r - 1); column B holds the
# compounding-growth formula that multiplies off the cell directly above
# it in column A:
#     B(r) = POWER(1.15, A(r-1)) * (100 * (A(r-1) * 1.01))
# Row 51 (A51/B51) is the last existing row, so we extend the same pattern
# down through row 101.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$firstNewRow = 52
$lastNewRow = 101
$lastExistingRow = 51

# 1) Copy the formatting (cell styles/number formats) of the last existing
#    row (A51:B51) down across all of the new rows in one shot, the same
#    way Excel would when you drag the fill handle down.
#    (NB: use ${var} inside "A$lastExistingRow:B..." style strings -- a bare
#    "$var:" is parsed as a PowerShell scope/drive qualifier, not the end of
#    the variable name.)
$ws.Range("A${lastExistingRow}:B${lastExistingRow}").Copy()
$ws.Range("A${firstNewRow}:B${lastNewRow}").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# 2) Fill in the actual index values (col A) and growth formulas (col B)
#    for each new row, referencing the cell directly above in col A --
#    matching the existing fill-down pattern used in rows 5-51.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $above = $r - 1
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 2).Formula = "=POWER(1.15,A$above)*(100*(A$above*1.01))"
}

# 3) Update the view: scroll so row 85 is at the top, and move the active
#    selection to C122 (matching where the author left off editing).
$win = $excel.ActiveWindow
$win.ScrollRow = 85
$win.ScrollColumn = 1
$null = $ws.Range("C122").Select()

Write-Output "Extended $($ws.Name) through row $lastNewRow"
